$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of this weekly block (before old row 798).
# This pushes the existing rows 798-864 down to 801-867, and the
# dimension grows from A1:T864 to A1:T867.
$ws.Rows("798:800").Insert()

# Fill in the newly inserted rows 798-800 with this week's new data.
# Common descriptive columns (A, B, C, E-J, Q, R, T) are the same for
# every row in this block.
$rows = @(798, 799, 800)
foreach ($r in $rows) {
    $ws.Range("A$r").Value = 11
    $ws.Range("B$r").Value = "Vega Monumental Concepción"
    $ws.Range("C$r").Value = "Bíobío"
    $ws.Range("D$r").Value = 45223
    $ws.Range("E$r").Value = 8
    $ws.Range("F$r").Value = "Fruta"
    $ws.Range("G$r").Value = 100108
    $ws.Range("H$r").Value = "Tropicales y subtropicales"
    $ws.Range("I$r").Value = 100108006
    $ws.Range("J$r").Value = "Plátano"
    $ws.Range("K$r").Value = "Sin especificar"
    $ws.Range("Q$r").Value = "`$/caja 20 kilos"
    $ws.Range("R$r").Value = "Ecuador"
    $ws.Range("T$r").Value = 20
}

$ws.Range("L798").Value = "Maduro"
$ws.Range("M798").Value = 100
$ws.Range("N798").Value = 20000
$ws.Range("O798").Value = 20000
$ws.Range("P798").Value = 20000
$ws.Range("S798").Value = 1000

$ws.Range("L799").Value = "Pintón"
$ws.Range("M799").Value = 400
$ws.Range("N799").Value = 21000
$ws.Range("O799").Value = 21000
$ws.Range("P799").Value = 21000
$ws.Range("S799").Value = 1050

$ws.Range("L800").Value = "Primera Pintón"
$ws.Range("M800").Value = 400
$ws.Range("N800").Value = 23000
$ws.Range("O800").Value = 23000
$ws.Range("P800").Value = 23000
$ws.Range("S800").Value = 1150
